# Add numero contrat generation function
# Refresh the "Etat Virement" sheet rows with the newly generated
# "N de contrat" values (and the other associated row data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "Karami abdelilah"
$ws.Range("B2").Value = "BB779645"
$ws.Range("C2").Clear()
$ws.Range("D2").Value = "chaabi"
$ws.Range("E2").Value = "chaabi"
$ws.Range("G2").Value = "001/TEST DR"
$ws.Range("I2").Value = 16000
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 15200

# --- Row 3 ---
$ws.Range("A3").Value = "mediexpets"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1196797"
$ws.Range("C3").Clear()
$ws.Range("D3").Value = "bmce"
$ws.Range("E3").Value = "bmce"
$ws.Range("G3").Value = "001/TEST DR"
$ws.Range("I3").Value = 4000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4000

# --- Row 4 ---
$ws.Range("A4").Value = "mediexpets"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "1196797"
$ws.Range("C4").Clear()
$ws.Range("D4").Value = "bmce"
$ws.Range("E4").Value = "bmce"
$ws.Range("F4").Value = "Logement de fonction"
$ws.Range("G4").Value = "001/LF/TEST DR"
$ws.Range("I4").Value = 120000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 120000

# --- Row 5 (totals) ---
$ws.Range("I5").Value = 140000
$ws.Range("J5").Value = 800
$ws.Range("K5").Value = 139200
